$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-05 Saturday" "2023-08-06 Sunday"

Replace-Text "91×75=" "39×35="
Replace-Text "55×78=" "72×69="
Replace-Text "52×89=" "81×36="
Replace-Text "64×76=" "70×60="
Replace-Text "28×15=" "57×27="

Replace-Text "95×33=" "66×27="
Replace-Text "82×94=" "65×93="
Replace-Text "55×97=" "23×66="
Replace-Text "15×83=" "59×26="
Replace-Text "68×60=" "53×62="

Replace-Text "25×72=" "50×93="
Replace-Text "76×62=" "50×36="
Replace-Text "34×84=" "58×77="
Replace-Text "13×78=" "25×22="
Replace-Text "80×31=" "11×27="

Replace-Text "39×99=" "51×63="
Replace-Text "88×85=" "47×56="
Replace-Text "53×60=" "31×50="
Replace-Text "46×19=" "19×65="
Replace-Text "46×89=" "49×31="

Replace-Text "87×87=" "99×23="
Replace-Text "64×91=" "38×88="
Replace-Text "87×69=" "16×59="
Replace-Text "82×75=" "42×67="
Replace-Text "37×36=" "56×12="
